# CS320-Sp22-102-roster.xlsx: "updated news, rosters, teams"
#
# Trey McBride / Bailey McBride was dropped from the section-102 roster.
# Removing the whole row shifts every row below it up by one, which is
# exactly what the author's diff shows (row 11 "McBride, Bailey" gone,
# rows 12-17 shift to 11-16, table shrinks from A1:E18 to A1:E17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 11 (the McBride row) the way a user would before deleting it,
# then delete the whole row - this shifts rows 12:17 up to 11:16 and leaves
# the selection sitting on the new row 11 (now "Mccloskey, Trey"), matching
# the saved <selection activeCell="A11" sqref="A11:XFD11"/>.
$null = $ws.Rows("11:11").Select()
$ws.Rows("11:11").Delete()

# The two mailto hyperlinks (originally on E9 "Kettula" and E16 "Stinson")
# need to track their cells. E9 is unaffected by the delete (it's above the
# removed row), but the Stinson link has to move from E16 down to E15 since
# everything below row 11 shifted up by one.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E9"), "mailto:jkettula@ycp.edu")
$ws.Hyperlinks.Add($ws.Range("E15"), "mailto:bstinson@ycp.edu")

# Re-adding hyperlinks stamps the built-in "Hyperlink" look (underline +
# theme color) on those cells; restore the original E-Mail-column
# formatting (bold, no underline, text number format) by copying the
# format from another cell in the same column that still has it.
$ws.Range("E4").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E15").PasteSpecial(-4122)
